# CCS-7, SP-30 Finalized proteomics format
# - Inserts a new "Reference Strain" property row into the metadata sheet
# - Renames the "Value Unit" example/description to drop "protein digest" -> "fmol/ug"
#   and tweaks the unit list description text (drop the stray "or" before AU)
# - Clears the now-removed "Scale" description (One of Lin, Log2, Log10, or Ln)
# - Leaves the data sheet text content unchanged (shared string indices shift
#   automatically as the strings in the pool are rearranged)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row 4 ("Reference Strain") above the current "Timepoint Type" row,
# copying the formatting from the row above (row 3, "Strain") so styles match.
$ws.Rows.Item(4).Insert()
$ws.Range("A3:C3").Copy()
$ws.Range("A4:C4").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Cells.Item(4, 1).Value = "Reference Strain"
$ws.Cells.Item(4, 2).Value = ""
$ws.Cells.Item(4, 3).Value = "The Reference Strain (for relative quantification data sets, leave empty for absolute)"

# Row 7 is now "Value Unit" (was row 6 before the insert)
$ws.Cells.Item(7, 2).Value = "fmol/ug"
$ws.Cells.Item(7, 3).Value = "One of mM, uM, Percent, RatioT1, RatioCs, AU, Dimensionless, fmol/ug"

# Row 8 is now "Scale"; its description column is cleared
$ws.Cells.Item(8, 3).Value = ""

# Match the new AutoFit-ish column width for column A
$ws.Columns.Item(1).ColumnWidth = 17.57421875

# Update the saved selection to match the author's cursor position
$ws.Activate()
$ws.Range("C16").Select()
